$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRow = 85

# Column A holds dates formatted as plain text (e.g. "2020-08-22") in the
# existing rows (shared strings, no explicit date number format). A direct
# $cell.Value = "2020-08-23" assignment gets auto-recognized as a real date
# by the COM layer's input heuristics and turned into a date serial number
# with a new number-format style. Routing the literal text through a scratch
# cell's Formula (so it's a computed string, not "typed" input) and then
# Copy / PasteSpecial (values only) into the target cell reproduces the
# original file's behavior: the cell keeps its default (no) style and its
# text is stored as a shared string, exactly like the existing rows above it.
$scratch = $ws.Cells.Item(200, 1)
$scratch.Formula = '="2020-08-23"'
$scratch.Copy()
$ws.Cells.Item($newRow, 1).PasteSpecial(-4163)
$scratch.Clear()

$ws.Cells.Item($newRow, 2).Value = 560164
$ws.Cells.Item($newRow, 3).Value = 618779
$ws.Cells.Item($newRow, 4).Value = 80198
$ws.Cells.Item($newRow, 5).Value = 60480
$ws.Cells.Item($newRow, 6).Value = 25.81
